# Apply the "23-11 15.13" -> "23-11 15.39" rerun results to the 4th sheet.
$wb = $excel.ActiveWorkbook

# Rename the sheet to reflect the new timestamp of the run.
$ws = $wb.Worksheets.Item("23-11 15.13")
$ws.Name = "23-11 15.39"

# Update the results produced for row 2 (Alternative model).
$ws.Range("D2").Value = 2
$ws.Range("F2").Value = 60.03342
$ws.Range("G2").Value = 0.14804
$ws.Range("H2").Value = 0.18421

# Update the results produced for row 3 (Standard model).
$ws.Range("D3").Value = 2
$ws.Range("F3").Value = 60.0083
$ws.Range("G3").Value = 0.24404
$ws.Range("H3").Value = 0.47368
$ws.Range("I3").Value = 7.13
